$d = $word.ActiveDocument

# Locate the "Docente(s) Responsável(eis)" heading paragraph and insert a
# new ListBullet paragraph right after it containing the professor entry.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Docente(s) Responsável(eis)*") {
        $target = $p
        break
    }
}

$target.Range.InsertParagraphAfter()
$newPara = $target.Next()
$newPara.Style = "ListBullet"
$newPara.Range.Text = "5111420 - Talita Martins Lacerda"
